$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# The "Name" column (A) for the three database rows was retyped to match
# the "User" column (C) values - finishing the mysql init database entries.
$ws.Range("A2").Value = "gameAdmin1"
$ws.Range("A3").Value = "gameAdmin2"
$ws.Range("A4").Value = "gameAdmin3"

# Leave the cursor on the last edited cell and make Database the active sheet.
$ws.Activate()
$ws.Range("A4").Select()

# Set up printing for the Database sheet (A4 paper, portrait) now that it's finished.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
